$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "wettened" -> "wetted"  (the run containing "wettene" loses its
#    trailing "ne"; formatting is untouched since it's the same run).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("wettene", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "wette", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "... work lightly & ..." -> "... work gently & ..."
#    The single run "lightly " (color=000000) is split into two runs:
#      - "gently"  (no explicit color)
#      - " "       (color=000000, unchanged)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("lightly") | Out-Null
$rng.Font.Color = -16777216   # wdColorAutomatic -> drop the explicit color
$rng.Text = "gently"

# ---------------------------------------------------------------------
# 3) "... or even better, white lead, melting ..." ->
#    "... or even better, lead white, melting ..."
#    "white lead" was split across two runs ("w" / "hite lead") that
#    already share identical formatting, so a plain replace merges them
#    into a single run with that same formatting.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("white lead", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "lead white", 2) | Out-Null
